$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) - update column F (想去人数)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1305
$ws.Range("F3").Value = 1187
$ws.Range("F4").Value = 891
$ws.Range("F5").Value = 108
$ws.Range("F7").Value = 660
$ws.Range("F11").Value = 2363
$ws.Range("F12").Value = 1593
$ws.Range("F13").Value = 1364
$ws.Range("F16").Value = 560
$ws.Range("F17").Value = 761
$ws.Range("F18").Value = 48
$ws.Range("F19").Value = 290
$ws.Range("F24").Value = 4701
$ws.Range("F25").Value = 213
$ws.Range("F26").Value = 197
$ws.Range("F27").Value = 50
$ws.Range("F30").Value = 209
$ws.Range("F31").Value = 87
$ws.Range("F33").Value = 682
$ws.Range("F34").Value = 1023
$ws.Range("F42").Value = 147

# Sheet "演出" (index 2) - update column F (想去人数)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 782
$ws.Range("F5").Value = 430
$ws.Range("F6").Value = 4

# Sheet "全部类型" (index 4) - update column F (想去人数)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1305
$ws.Range("F4").Value = 782
$ws.Range("F5").Value = 1187
$ws.Range("F6").Value = 891
$ws.Range("F8").Value = 430
$ws.Range("F9").Value = 108
$ws.Range("F11").Value = 660
$ws.Range("F14").Value = 4
$ws.Range("F18").Value = 2363
$ws.Range("F19").Value = 1593
$ws.Range("F20").Value = 1364
$ws.Range("F23").Value = 560
$ws.Range("F25").Value = 761
$ws.Range("F26").Value = 48
$ws.Range("F27").Value = 290
$ws.Range("F30").Value = 4701
$ws.Range("F31").Value = 213
$ws.Range("F32").Value = 50
$ws.Range("F35").Value = 209
$ws.Range("F36").Value = 87
$ws.Range("F38").Value = 682
$ws.Range("F45").Value = 147
